$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new week of price data (2022-02-11, serial 44603) was recorded for this
# product/market combination. It is inserted right above the existing block
# of rows for this series (the sheet is otherwise sorted by date), so two new
# rows are inserted at row 214 and every subsequent row shifts down by two
# (old row 331 -> new row 333; dimension grows from A1:R331 to A1:R333).
$ws.Rows.Item(214).Resize(2).Insert()

# Row 214: "Primera" quality for the new week
$ws.Range('A214').Value = 8
$ws.Range('B214').Value = 'Terminal La Palmera de La Serena'
$ws.Range('C214').Value = 'Coquimbo'
$ws.Range('D214').Value = 44603
$ws.Range('E214').Value = 4
$ws.Range('F214').Value = 100112009
$ws.Range('G214').Value = 'Acelga'
$ws.Range('H214').Value = 'Sin especificar'
$ws.Range('I214').Value = 'Primera'
$ws.Range('J214').Value = 2400
$ws.Range('K214').Value = 500
$ws.Range('L214').Value = 600
$ws.Range('M214').Value = 550
$ws.Range('N214').Value = '$/atado 1,5 a 2 kilos'
$ws.Range('O214').Value = 'Provincia del Elquí'
$ws.Range('P214').Value = 275
$ws.Range('Q214').Value = 2
$ws.Range('R214').Value = 'Hortaliza'

# Row 215: "Segunda" quality for the new week
$ws.Range('A215').Value = 8
$ws.Range('B215').Value = 'Terminal La Palmera de La Serena'
$ws.Range('C215').Value = 'Coquimbo'
$ws.Range('D215').Value = 44603
$ws.Range('E215').Value = 4
$ws.Range('F215').Value = 100112009
$ws.Range('G215').Value = 'Acelga'
$ws.Range('H215').Value = 'Sin especificar'
$ws.Range('I215').Value = 'Segunda'
$ws.Range('J215').Value = 1340
$ws.Range('K215').Value = 400
$ws.Range('L215').Value = 450
$ws.Range('M215').Value = 425
$ws.Range('N215').Value = '$/atado 1,5 a 2 kilos'
$ws.Range('O215').Value = 'Provincia del Elquí'
$ws.Range('P215').Value = 212
$ws.Range('Q215').Value = 2
$ws.Range('R215').Value = 'Hortaliza'
